$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 + season-sim values to OFF/DEF run & pass logs ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '4 2 2 1 -3 7 3 5 6 0 0 8 10 2 5 4 -1 1 1 1 12 5 4 1 3 1 8 9 3 14 9 2 1 0 11 0 0 4 5 3 9 2 5 2 7 9 5 6 7 2 4 2 0 -2 13 29 1 7 0 2 3 5 -5 6 1 -1 1 1 0 14 4 3 0 1 1 1 1 15 3 2 4 4 0 7 3 5 2 11 7 4 11 -1 -1 3 10 9 5 7 -2 2 10 -1 6 3 4 1 -2 2 17 5 -1 4 4 3 0 12 2 5 4 21 2 12 1 1 0 0 2 9 0 0 5 -2 15 11 11 5 1 9 1 1 4 5 2 7 2 -1 4 5 9 -1 3 4 1 4 5 3 4 3 2 3 4 6 6 1 0 0 3 -5 4 5 1 1 0 3 6 2 1 6 0 2 -2 4 2 5 2 7 4 19 3 8 -1 -2 3 4 -2 3 1 6 17 8 1 1 1 2 3 18 1 8 2 0 6 3 0 6 3 0 2 6 18 0 1 11 2 5 3 -1 3 17 -1 -3 -1 4 4 7 2 1 4 3 -3 3 2 14 2 0 4 3 1 5 4 4 3 3 4 -2 12 2 5 -3 9 3 16 8 4 1 2 4 12 3 1 8 5 3 3 4 2 4 11 4 0 1 1 3 5 0 4 0 0 14 1 4 2 6 1 26 7 1 4 1 7 7 1 -2 6 2 2 2 2 2 3 5 2 3 7 2 9 3 5 2 1 9 1 5 6 5 4 6 7 2 10 4 8 4 6 3 4 5 4 9 14 2 7 5 1 31 8 3 2 0 12 3 4 3 13 24 9 4 2 1 0 1 2 6 -4 5 0 14 8 3 -1 8 8 9 -1 -1 24 1 2 1 7 8 2 4 5 0 3 5 1 -1 4 0 3 2 -7 5 1 7 16 1 3 2 2 16 6 2 15 4 -1 8 3 -1 -1 2 2 1 15 8 4 -3 5 2 5 4 2 6 18 -1 2 -1 6 2 1 3 2 8 3 4 -2 6 1 12 1 6 5 2 2 3 24 4 1 6 0 1 14 10 1 3 5 3 19 4 4 0 -1 0 1 15 4 2 4 2 0 6 3 -2 4 -1 4 3 2 0 0 4 2 4 0 3 1 13 7 4 5 -1 3 7 4 5 2 16 4 3 4 1 8 1 7 5 0 1 3 4 4 3 3 6 5 2 8 3 7 5 2 2 3 8 6 0 16 3 6 2 3 9 5 3 9 3 10 -5 4 1 7 5 5 5 -3 2 2 5 3 -2 12 3 1 1 0 2 6 2 0 3 1 7 0 1 2 4 4 6 -2 6 0 5 2 0 -2 7 2 2 2 5 -5 0 4 3 0 -1 2 3 1 -1 5 4 2 4 2 1 -1 4 1 14 0 11 0 1 7 2 1 14 3 -1 0 1 -3 5 0 4 1 2 7 3 2 20 9 9 0 -3 6 4 0 2 2 4 3 4 6 3 3 6 0 0 3 2 6 3 1 5 0 2 2 3 2 5 -3 3 0 2 16 0 2 4 4 3 3 8 7 2 0 8 2 1 9 2 1 1 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2 3 4 0 0 4 10 1 1 12 4 5 0 8 2 5 4 1 6 6 9 3 7 2 0 11 3 2 0 3 24 8 30 4 12 1 4 2 2 1 4 2 2 4 3 4 4 2 9 4 -3 1 10 6 2 0 4 4 1 5 3 4 0 0 5 -1 2'
$wsYDS.Range("C2").Value = '7 10 5 4 2 5 4 2 15 -1 7 7 4 3 1 4 5 -2 6 5 6 7 11 1 8 5 2 11 2 2 23 7 4 12 5 -1 5 1 -1 11 14 3 1 16 0 1 3 14 6 4 18 3 4 2 4 3 0 2 10 2 3 4 10 3 3 16 2 11 1 4 2 3 -1 5 2 2 4 1 7 2 1 -1 2 6 7 8 0 4 8 9 2 3 18 -4 8 8 2 1 4 3 3 5 0 7 5 5 37 5 7 11 4 0 9 4 3 7 9 9 -1 5 3 5 1 3 0 -2 5 -1 4 3 6 -2 12 1 8 34 1 14 8 -3 0 4 4 4 7 4 18 6 5 12 0 -2 4 10 5 5 4 2 11 0 -3 3 4 2 3 4 5 2 4 4 9 14 3 4 1 1 10 5 3 4 3 4 3 16 6 -6 4 1 5 4 4 5 -3 3 12 4 2 -1 8 4 2 3 6 3 3 2 7 4 3 4 -2 7 2 10 3 1 3 4 0 3 8 -1 7 4 7 4 4 -2 2 11 0 1 20 -2 14 9 1 2 3 -2 5 25 7 5 18 3 2 3 -1 9 4 2 20 1 4 9 20 -2 8 0 3 -2 3 5 6 1 5 2 1 3 3 4 1 1 4 14 17 0 1 2 10 1 2 2 9 1 -1 0 3 2 5 1 3 3 0 10 0 -1 -1 4 2 6 3 13 32 0 8 5 9 -2 3 -1 9 2 12 4 0 0 0 -3 -8 1 -1 6 16 1 4 17 4 3 8 -1 4 2 3 4 2 4 11 10 2 13 2 3 3 13 2 1 7 0 15 10 18 3 1 9 0 1 3 3 2 7 2 2 5 5 -1 4 2 3 3 2 4 18 0 -3 1 0 7 18 3 2 1 1 7 5 35 3 6 1 2 2 11 4 5 2 0 3 4 4 1 4 5 1 -1 10 -4 1 3 3 1 8 1 3 5 2 46 4 6 1 10 1 -1 2 5 4 0 8 4 4 4 1 4 15 2 5 0 0 7 3 4 2 8 1 3 2 3 5 3 0 4 2 6 4 2 5 2 14 7 4 1 19 3 2 0 1 4 -2 1 6 0 1 27 8 1 4 4 8 0 2 10 1 1 23 38 1 6 9 -1 0 12 -5 6 2 4 7 5 0 -5 6 0 -3 -1 1 10 1 17 4 8 4 4 4 6 4 -1 13 13 3 6 1 4 5 5 7 0 5 3 1 2 4 24 1 4 9 3 12 1 2 2 5 -1 1 7 3 13 1 2 7 2 2 2 2 1 3 11 2 2 3 5 3 2 4 -1 6 34 1 0 -1 -2 6 5 7 13 2 6 2 12 3 -1 4 0 7 4 8 3 9 3 3 1 1 10 2 4 8 1 0 2 4 1 9 4 19 4 12 2 8 4 2 -2 0 0 5 2 11 0 -4 1 3 1 11 5 7 7 5 39 -1 5 4 15 1 1 4 2 5 0 2 5 3 4 1 4 4 4 2 1 15 5 2 2 3 1 1 -3 0 5 -1 5 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3 6 8 11 1 2 2 3 1 5 3 5 7 4 1 -5 11 0 4 6 6 20 -2 2 1 3 4 3 1 -1 7 7 4 4 0 2 2 10 3 1 2 19 0 2 1 2 4'
$wsYDS.Range("B3").Value = '3 12 6 13 14 -1 18 20 15 -1 9 7 3 14 5 3 18 8 4 21 17 2 13 13 26 2 4 14 2 -1 7 27 15 15 15 12 20 14 23 9 0 12 11 3 10 8 8 19 3 13 2 12 8 7 20 19 3 11 8 17 11 4 8 15 4 3 3 16 8 3 0 9 5 9 18 15 15 26 11 17 9 3 12 21 6 5 12 1 4 17 4 9 -4 16 9 11 10 18 18 8 47 15 3 28 6 6 8 12 19 2 22 70 4 1 8 4 7 6 15 31 32 4 6 8 6 4 10 10 16 3 -3 3 43 4 7 30 15 18 4 7 2 7 8 0 8 15 3 9 7 10 7 8 10 8 6 5 11 17 19 35 3 9 4 16 16 17 9 14 5 17 6 11 19 9 23 13 15 3 23 8 7 13 17 8 19 12 5 2 1 3 13 8 4 9 3 6 10 8 13 6 7 0 25 11 15 1 8 6 9 21 12 2 2 12 12 8 5 22 5 11 7 6 5 19 17 13 19 13 4 25 10 7 7 6 15 7 4 8 6 15 8 16 20 7 9 2 11 5 5 17 14 35 5 22 11 8 20 4 8 9 9 18 16 14 18 7 7 10 14 2 12 16 9 10 10 -1 11 11 13 5 23 10 29 2 6 9 9 18 4 22 3 7 8 15 3 15 3 14 5 10 5 8 4 8 6 1 10 4 7 1 1 13 5 5 6 2 3 5 4 5 6 14 5 8 10 1 4 13 8 17 31 7 9 59 34 11 6 8 7 22 3 3 1 10 8 5 32 5 17 13 25 5 6 2 16 11 31 7 12 13 23 7 6 5 -1 7 10 7 9 18 25 17 5 9 18 23 6 36 4 0 30 12 3 -2 15 3 13 13 7 6 9 5 18 12 4 5 14 4 11 16 11 3 5 8 0 4 2 7 8 4 5 6 8 5 -1 12 8 4 10 9 7 4 13 4 4 2 7 8 9 2 2 3 5 9 6 7 4 10 3 7 1 5 2 18 4 10 8 6 4 -7 11 8 7 4 22 25 42 22 4 21 4 34 9 12 10 5 24 5 3 12 5 2 21 6 10 8 4 23 10 2 7 1 0 11 5 21 7 10 17 4 20 6 6 6 22 25 15 3 2 8 1 5 20 7 8 5 14 0 6 32 20 -1 12 18 2 12 9 8 5 -1 9 17 5 2 10 10 0 7 2 15 9 7 19 5 8 11 5 14 6 6 9 18 27 9 11 1 7 10 17 16 4 23 7 8 0 14 15 22 12 8 3 8 5 1 -1 6 12 10 4 10 -2 8 17 18 12 3 6 12 10 5 11 18 21 5 6 10 9 1 4 9 11 2 13 14 6 6 7 16 9 11 7 21 2 6 19 2 52 11 9 1 6 14 35 23 6 64 3 26 1 19 7 5 8 13 3 4 7 7 6 9 1 8 65 15 8 0 8 1 7 6 5 7 6 2 5 -1 25 4 10 6 0 9 57 9 5 5 5 7 11 11 7 7 4 11 11 9 2 -2 5 6 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17 5 1 9 37 15 12 4 27 20 20 7 1 17 9 11 1 4 4 5 -1 18 7 11 6 7 6 10 9 17 13 40 24 1 5 12'
$wsYDS.Range("C3").Value = '10 25 6 6 9 16 10 16 11 9 13 4 5 6 9 5 25 15 26 8 5 1 36 4 11 46 22 16 15 12 9 8 6 47 24 6 21 1 46 24 2 16 3 11 9 6 6 6 7 3 10 7 11 16 16 8 5 6 12 9 5 4 8 6 6 16 9 18 5 37 5 23 14 4 19 18 1 4 9 11 6 11 57 4 3 21 11 9 30 9 17 32 15 5 11 15 11 20 0 12 19 15 12 5 9 9 4 9 1 9 18 2 9 8 1 2 9 4 8 13 6 22 6 3 12 21 8 7 1 17 14 8 11 22 8 1 9 13 9 5 8 3 7 13 6 26 34 7 6 8 15 15 19 11 9 13 9 7 15 5 11 14 4 6 4 8 10 17 4 16 56 16 10 4 6 17 7 14 21 11 6 4 10 9 7 35 6 11 17 9 5 5 6 15 9 10 9 16 -2 7 4 4 28 13 13 15 2 9 15 4 41 9 22 1 15 9 6 9 16 7 0 29 16 37 6 6 19 0 30 6 9 27 12 3 3 1 15 7 8 3 12 7 4 5 72 15 11 7 9 8 -3 5 6 2 1 16 17 7 37 21 9 13 16 9 5 15 6 26 44 28 26 5 32 13 5 11 9 -2 11 3 8 5 8 21 12 20 35 8 23 26 8 11 4 12 -1 16 11 4 29 2 19 12 19 9 9 13 22 11 33 10 3 22 4 85 4 7 7 5 13 5 6 8 18 7 8 10 16 19 6 14 16 27 32 15 19 3 51 20 56 4 7 4 9 8 5 22 6 7 12 8 25 9 7 21 26 9 6 17 16 4 1 6 10 3 8 3 7 11 1 35 5 6 1 7 7 22 13 6 7 8 6 9 2 41 3 24 7 23 6 10 25 8 18 8 5 23 12 21 1 4 15 10 12 23 7 23 16 3 16 32 34 1 7 6 12 1 5 4 28 6 3 2 3 17 8 -2 3 14 10 1 3 41 11 24 20 9 26 13 10 10 18 6 15 62 6 10 10 11 12 3 4 11 5 12 6 14 12 16 4 15 7 34 19 9 22 24 1 19 4 6 8 9 11 20 10 28 11 29 11 11 21 7 14 10 7 10 20 12 -3 23 11 7 2 13 12 6 13 15 5 39 9 6 4 4 26 49 4 11 12 6 11 23 28 1 11 8 1 20 3 15 7 8 6 15 4 10 6 3 3 8 8 17 11 14 19 13 1 11 -1 20 8 4 3 23 9 9 15 6 7 5 10 6 9 8 15 4 39 15 4 10 3 5 10 7 4 7 4 11 15 7 5 6 13 3 11 9 9 11 0 8 7 8 4 8 20 5 8 30 6 9 15 9 2 5 5 29 14 20 9 5 22 5 4 9 23 5 5 2 62 12 8 11 8 15 64 4 13 4 7 18 24 13 17 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11 16 2 14 19 4 11 6 23 1 9 24 10 4 17 8 11 12 4 14 -5 6 4 4 56'

# --- ST sheet: append Week 16 + season-sim values to RA/RM logs ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value = '17 64 60 63 71 66 62 67 53 64 74 61 65 66 46 66 66 67 49 68 58 41 70 58 60 62 61 66 67 57 66 66 62 63 57 66 69 40 66 67 71'
$wsST.Range("B5").Value = '-1 24 6 29 20 25 0 21 15 13 22 14 31 26 6 15 35 24 15 27 5 17 23 15 21 23 0 23 27 5 24 21 19 18 18 21 29 20 20 31 26'
$wsST.Range("B6").Value = '19 12 5 15 23 4 15 22 22 21 13 11 31 15 18 20 13 18 7 22 25 16 15 17 23'
$wsST.Range("D3").Value = '45 58 49 63 50 51 58 40 38 42 37 47 51 53 31 44 36 29 45 45 37 46 32 63 38 44 40 47 40 39 44 41 40 62 42 52 43 51 42 54 42 34 43 34 31 55 52 41 32 54 42 37 50 51 45 54 41 51 46 46 46 44 32 42 34 55 55 40 52 37 51 48 47 30 36 61 39 34 50 32 41 46 53 56 39 35 36 44 30 37 32 59 49 47 58 43 59 58 55 44 49 39 33 50 36 38 40 38 49 44 40 52 47 53 41 57 40 37 44 35 46 53 50 65 39 37 46 51 52 56 41 37 47 43 42 45'
$wsST.Range("D4").Value = '0 0 0 0 7 0 11 0 0 0 8 0 0 17 0 0 0 0 0 0 0 0 0 0 0 0 12 2 0 0 0 -4 0 0 0 1 0 9 0 0 0 0 0 0 0 9 4 0 0 0 0 67 0 0 3 13 0 32 -3 13 9 19 0 0 0 84 0 7 0 -1 14 20 0 0 0 13 7 0 2 0 0 0 0 -1 0 0 0 0 0 0 0 0 0 17 0 16 0 9 0 0 0 2 0 0 0 0 0 0 0 9 0 0 0 9 0 19 0 0 0 0 0 0 15 11 0 -4 0 0 11 12 10 0 0 0 0 17'
$wsST.Range("D5").Value = '0 0 -3 0 1 0 10 0 13 15 0 11 6 5 16 1 0 0 30 0 14 15 7 0 0 0 0 88 0 0 10 0 18 20 19 0 6 0 0 0 12 7 0 0 2 14 6 0 0 0 -4 11 1 0 0 7 0 0 0 0 0 0 12 0 18 4 0 0 8 0 7 12 -1 14 0 0 0 0 0 0 0 16 0 8 0 -4 8 15 0 10 18 0 6 0 0 0 0 0 4 0 0 0 0 0 7 7 0 0 0 0 7 0 0 0 0 0 0 0 0 5 -1 9 0 0 0 0 10'

# --- OFF sheet: updated season totals after Week 16 + season sim ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 352
$wsOFF.Range("D2").Value = 25
$wsOFF.Range("F2").Value = 124
$wsOFF.Range("G2").Value = 109
$wsOFF.Range("H2").Value = 10
$wsOFF.Range("I2").Value = 15
$wsOFF.Range("J2").Value = 64
$wsOFF.Range("N2").Value = 32
$wsOFF.Range("B3").Value = 20
$wsOFF.Range("C3").Value = 386
$wsOFF.Range("E3").Value = 48
$wsOFF.Range("F3").Value = 206
$wsOFF.Range("G3").Value = 90
$wsOFF.Range("I3").Value = 102
$wsOFF.Range("J3").Value = 129
$wsOFF.Range("L3").Value = 538
$wsOFF.Range("M3").Value = 364
$wsOFF.Range("Q3").Value = 931

# --- DEF sheet: updated season totals after Week 16 + season sim ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 371
$wsDEF.Range("F2").Value = 107
$wsDEF.Range("G2").Value = 89
$wsDEF.Range("H2").Value = 10
$wsDEF.Range("J2").Value = 41
$wsDEF.Range("N2").Value = 37
$wsDEF.Range("O2").Value = 45
$wsDEF.Range("P2").Value = 22
$wsDEF.Range("C3").Value = 373
$wsDEF.Range("E3").Value = 70
$wsDEF.Range("F3").Value = 224
$wsDEF.Range("H3").Value = 63
$wsDEF.Range("I3").Value = 117
$wsDEF.Range("J3").Value = 98
$wsDEF.Range("L3").Value = 522
$wsDEF.Range("M3").Value = 338
$wsDEF.Range("Q3").Value = 937

# --- ST sheet: updated season totals after Week 16 + season sim ---
$wsST.Range("B2").Value = 155
$wsST.Range("D2").Value = 136
$wsST.Range("F2").Value = 134
$wsST.Range("G2").Value = 131
$wsST.Range("J2").Value = 55
$wsST.Range("K2").Value = 52
$wsST.Range("L2").Value = 41
$wsST.Range("M2").Value = 31
$wsST.Range("N2").Value = 20
$wsST.Range("B3").Value = 114

# --- TURNS sheet: updated season totals after Week 16 + season sim ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 14
$wsTURNS.Range("C3").Value = 10
$wsTURNS.Range("D3").Value = 17
$wsTURNS.Range("E3").Value = 16

# --- PEN sheet: updated season totals after Week 16 + season sim ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B3").Value = 31

